# feat: add parameters to check box tests
# Add a new "CheckBox" worksheet at the end of the workbook, populate it
# with header + parameter rows, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last sheet so it lands at the end.
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "CheckBox"

# Header row (bold).
$ws.Range("A1").Value = "Box"
$ws.Range("B1").Value = "Results"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows -- written in the same order the original authoring session
# used (matches shared-string table ordering in the target workbook).
$ws.Range("A4").Value = "Desktop React"
$ws.Range("B4").Value = "desktop notes commands react"

$ws.Range("A5").Value = "Classified"
$ws.Range("B5").Value = "classified"

$ws.Range("B2").Value = "home desktop notes commands documents workspace react angular veu office public private classified general downloads wordFile excelFile"
$ws.Range("A2").Value = "Home"

$ws.Range("A3").Value = "Home WorkSpace Public"
$ws.Range("B3").Value = "desktop notes commands private classified general downloads wordFile excelFile"

# Column widths to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 13.71
$ws.Columns.Item(2).ColumnWidth = 29

# Page setup mirrors the other sheets in the workbook (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Select a harmless cell and make this sheet the active tab (mirrors the
# "tabSelected" / activeTab move away from userPositive onto CheckBox).
$ws.Range("A7").Select()
$ws.Activate()
